$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 1270
$ws1.Range("F4").Value = 2763

# Sheet "演出" (Performance)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("G2").Value = "不可售"

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("G3").Value = "不可售"
$ws4.Range("F5").Value = 1270
$ws4.Range("F6").Value = 2763
